$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Sector value ("Software") for rows 2 and 3 in column J
$ws.Range("J2").Value = "Software"
$ws.Range("J3").Value = "Software"

# Update the selected cell/range to match the authored change
$ws.Range("J3").Select()
